$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (id) and C (speaker_variant) for rows 2..12.
# Column D (is_prefered) is cleared for every one of these rows.
$data = @(
    @{ Row = 2;  B = "#rykert";   C = "Rykert" },
    @{ Row = 3;  B = "#krispyn";  C = "Krispyn" },
    @{ Row = 4;  B = "#lyseeth";  C = "Lyseeth" },
    @{ Row = 5;  B = "#valerius"; C = "Valerius" },
    @{ Row = 6;  B = "#klareth";  C = "Klareth" },
    @{ Row = 7;  B = "#lysbeth";  C = "Lysbeth" },
    @{ Row = 8;  B = "#klaret";   C = "Klaret" },
    @{ Row = 9;  B = "#gustaaf";  C = "Gustaaf" },
    @{ Row = 10; B = "#sofy";     C = "Sofy" },
    @{ Row = 11; B = "#sofy";     C = "sofy" },
    @{ Row = 12; B = "#rudolf";   C = "Rudolf" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = ""
}
